# Add I0 (I) and IF (J) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy formats from an existing header cell (H1) so the new
# header cells share the same style as the others, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J47
$iValues = @(6,3,8,7,6,2,6,8,5,6,8,7,6,7,9,6,7,8,7,5,7,13,8,8,6,8,6,7,1,8,7,7,4,5,3,4,9,4,10,8,4,4,4,5,4,4)
$jValues = @(7,4,8,8,7,5,7,8,6,6,8,8,7,9,9,7,7,8,8,6,8,13,9,9,8,9,8,9,2,8,7,8,6,7,4,6,9,7,10,9,5,5,5,6,5,5)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
